$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix typo in the terminal string literal for "ifj21"
$ws.Range("C2").Value = '"ifj21"'

# 2. Re-arrange the <stmt> / new <ID_assign_or_fun> productions (rows 27-32).
#    Row 27 used to be "ID <ID_next> = <expr> <expr2>"; it now becomes the
#    "if" production. Row 28 ("if...") becomes the "while" production.
#    Row 29 ("while...") becomes "ID <ID_assign_or_fun>". Row 30
#    ("ID <ID_next> = <fun_ call>") becomes "return <expr> <expr2>".
#    Row 31 ("<fun_ call>") becomes the new nonterminal's first alt.
#    Row 32 ("return...") becomes the new nonterminal's second alt.

# Row 27: <stmt> -> if term then <stmts> else <stmts> end
$ws.Cells.Item(27,1).Value = "<stmt>"
$ws.Cells.Item(27,2).Value = "if"
$ws.Cells.Item(27,3).Value = "term"
$ws.Cells.Item(27,4).Value = "then"
$ws.Cells.Item(27,5).Value = "<stmts>"
$ws.Cells.Item(27,6).Value = "else"
$ws.Cells.Item(27,7).Value = "<stmts>"
$ws.Cells.Item(27,8).Value = "end"

# Row 28: <stmt> -> while term do <stmts> end
$ws.Cells.Item(28,1).Value = "<stmt>"
$ws.Cells.Item(28,2).Value = "while"
$ws.Cells.Item(28,3).Value = "term"
$ws.Cells.Item(28,4).Value = "do"
$ws.Cells.Item(28,5).Value = "<stmts>"
$ws.Cells.Item(28,6).Value = "end"
$ws.Cells.Item(28,7).Value = $null
$ws.Cells.Item(28,8).Value = $null

# Row 29: <stmt> -> ID <ID_assign_or_fun>
$ws.Cells.Item(29,1).Value = "<stmt>"
$ws.Cells.Item(29,2).Value = "ID"
$ws.Cells.Item(29,3).Value = "<ID_assign_or_fun>"
$ws.Cells.Item(29,4).Value = $null
$ws.Cells.Item(29,5).Value = $null
$ws.Cells.Item(29,6).Value = $null

# Row 30: <stmt> -> return <expr> <expr2>
$ws.Cells.Item(30,1).Value = "<stmt>"
$ws.Cells.Item(30,2).Value = "return"
$ws.Cells.Item(30,3).Value = "<expr>"
$ws.Cells.Item(30,4).Value = "<expr2>"
$ws.Cells.Item(30,5).Value = $null

# Row 31: <ID_assign_or_fun> -> <fun_ call>
$ws.Cells.Item(31,1).Value = "<ID_assign_or_fun>"
$ws.Cells.Item(31,2).Value = "<fun_ call>"

# Row 32: <ID_assign_or_fun> -> <ID_next> = <expr> <expr2>
$ws.Cells.Item(32,1).Value = "<ID_assign_or_fun>"
$ws.Cells.Item(32,2).Value = "<ID_next>"
$ws.Cells.Item(32,3).Value = "="
$ws.Cells.Item(32,4).Value = "<expr>"
$ws.Cells.Item(32,5).Value = "<expr2>"

# 3. Fix the mis-named terminal "a" -> "term" in the <expr> production
#    (currently row 38, column B).
$ws.Cells.Item(38,2).Value = "term"

# 4. Remove the stray blank row 33, shifting rows 34+ up by one.
$ws.Rows(33).EntireRow.Delete()

# 5. Update the saved selection/view.
$ws.Range("C2").Select()

# 6. Page setup (paper size / orientation) now present on the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
